$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "31.400.04"
$ws.Range("E2").Value = "  +3.67%  "
# Row 3
$ws.Range("D3").Value = "2.006.36"
$ws.Range("E3").Value = "  +7.47%  "
# Row 4
$ws.Range("E4").Value = "  -0.28%  "
# Row 5
$ws.Range("D5").Value = "'0.7863"
$ws.Range("E5").Value = "  +66.47%  "
# Row 6
$ws.Range("D6").Value = "'259.88"
$ws.Range("E6").Value = "  +6.90%  "
# Row 7
$ws.Range("D7").Value = "'0.9974"
$ws.Range("E7").Value = "  -0.32%  "
# Row 8
$ws.Range("D8").Value = "'0.3635"
$ws.Range("E8").Value = "  +26.74%  "
# Row 9
$ws.Range("D9").Value = "'28.41"
$ws.Range("E9").Value = "  +32.23%  "
# Row 10
$ws.Range("D10").Value = "'0.07069"
$ws.Range("E10").Value = "  +9.31%  "
# Row 11
$ws.Range("D11").Value = "'0.8410"
$ws.Range("E11").Value = "  +17.50%  "
# Row 12
$ws.Range("D12").Value = "'0.08070"
$ws.Range("E12").Value = "  +3.68%  "
# Row 13
$ws.Range("D13").Value = "2.003.99"
$ws.Range("E13").Value = "  +7.39%  "
# Row 14
$ws.Range("D14").Value = "'101.12"
$ws.Range("E14").Value = "  +5.08%  "
# Row 15
$ws.Range("D15").Value = "'5.669"
$ws.Range("E15").Value = "  +10.94%  "
# Row 16
$ws.Range("D16").Value = "'274.10"
$ws.Range("E16").Value = "  -2.38%  "
# Row 17
$ws.Range("D17").Value = "31.385.88"
$ws.Range("E17").Value = "  +3.66%  "
# Row 18
$ws.Range("D18").Value = "'14.75"
$ws.Range("E18").Value = "  +14.03%  "
# Row 19
$ws.Range("D19").Value = "'5.936"
$ws.Range("E19").Value = "  +13.42%  "
# Row 20
$ws.Range("D20").Value = "'0.000007963"
$ws.Range("E20").Value = "  +6.87%  "
# Row 21
$ws.Range("D21").Value = "2.265.46"
$ws.Range("E21").Value = "  +7.56%  "
# Row 22
$ws.Range("D22").Value = "'0.9977"
$ws.Range("E22").Value = "  -0.27%  "
# Row 23
$ws.Range("D23").Value = "'0.9980"
$ws.Range("E23").Value = "  -0.29%  "
# Row 24
$ws.Range("D24").Value = "'7.211"
$ws.Range("E24").Value = "  +15.72%  "
# Row 25
$ws.Range("D25").Value = "'10.22"
$ws.Range("E25").Value = "  +14.14%  "
# Row 26
$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D26").Value = "'0.1505"
$ws.Range("E26").Value = "  +56.93%  "
# Row 27
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").Value = "'164.01"
$ws.Range("E27").Value = "  +1.24%  "
# Row 28
$ws.Range("D28").Value = "'20.08"
$ws.Range("E28").Value = "  +7.70%  "
# Row 29
$ws.Range("D29").Value = "'2.411"
$ws.Range("E29").Value = "  +29.01%  "
# Row 30
$ws.Range("D30").Value = "'1.629"
# Row 31
$ws.Range("D31").Value = "'4.630"
$ws.Range("E31").Value = "  +10.40%  "
# Row 32
$ws.Range("E32").Value = "  +3.18%  "
# Row 33
$ws.Range("D33").Value = "'4.431"
$ws.Range("E33").Value = "  +7.89%  "
# Row 34
$ws.Range("D34").Value = "'0.05216"
$ws.Range("E34").Value = "  +9.32%  "
# Row 35
$ws.Range("D35").Value = "'1.227"
$ws.Range("E35").Value = "  +10.14%  "
# Row 36
$ws.Range("D36").Value = "'0.7642"
$ws.Range("E36").Value = "  +12.29%  "
# Row 37
$ws.Range("D37").Value = "'2.810"
$ws.Range("E37").Value = "  +3.61%  "
# Row 38
$ws.Range("D38").Value = "'0.02016"
$ws.Range("E38").Value = "  +7.41%  "
# Row 39
$ws.Range("D39").Value = "'2.958"
$ws.Range("E39").Value = "  +4.40%  "
# Row 40
$ws.Range("D40").Value = "'6.714"
$ws.Range("E40").Value = "  +8.53%  "
# Row 41
$ws.Range("D41").Value = "'80.56"
$ws.Range("E41").Value = "  +7.57%  "
# Row 42
$ws.Range("D42").Value = "'2.187"
$ws.Range("E42").Value = "  +14.29%  "
# Row 43
$ws.Range("D43").Value = "'0.4742"
$ws.Range("E43").Value = "  +13.47%  "
# Row 44
$ws.Range("D44").Value = "'0.8636"
$ws.Range("E44").Value = "  +5.03%  "
# Row 45
$ws.Range("D45").Value = "'104.66"
$ws.Range("E45").Value = "  +4.48%  "
# Row 46
$ws.Range("D46").Value = "'0.9983"
$ws.Range("E46").Value = "  -0.12%  "
# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.04"
$ws.Range("E47").Value = "  +5.00%  "
# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.690"
$ws.Range("E48").Value = "  +10.64%  "
# Row 49
$ws.Range("D49").Value = "'0.4372"
$ws.Range("E49").Value = "  +13.31%  "
# Row 50
$ws.Range("D50").Value = "'37.11"
$ws.Range("E50").Value = "  +6.37%  "
# Row 51
$ws.Range("D51").Value = "'942.04"
$ws.Range("E51").Value = "  +7.13%  "
